# Update the "Förändrad" (Changed) date column (C) for every data row
# (the header is row 1) from 2026-02-08 (serial 46061) to 2026-02-09
# (serial 46062).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 46062
}
